# Build site at 2022-09-26 16:07:08 UTC
# Re-arranges rows 10-25 of the LOQ4086 sheet: the "Docentes responsaveis" /
# "Programa" / "Avaliacao" blocks are reshuffled and two trailing rows are
# removed (final dimension becomes A1:C23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Remove the two trailing rows (old rows 24 and 25) completely, so the
#    sheet shrinks from A1:C25 down to A1:C23, matching the new dimension.
# ---------------------------------------------------------------------
$ws.Range("A24:C25").EntireRow.Delete()

# ---------------------------------------------------------------------
# 2) Row 10 / C10: objective text is replaced by the professor entry.
# ---------------------------------------------------------------------
$teacher1 = "8151869 - Livia Chaguri e Carvalho"
$ws.Range("B10").Value = $teacher1
$ws.Range("C10").Value = $teacher1

# ---------------------------------------------------------------------
# 3) Rows 13-23 are rewritten entirely (labels in column A, values in
#    columns B/C) to match the new layout, including row heights.
# ---------------------------------------------------------------------

# Row 13: "Programa resumido:" label, value becomes the activation date.
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "01/01/2016"
$ws.Range("C13").Value = "01/01/2016"
$ws.Rows.Item(13).RowHeight = 60

# Row 14: "Short syllabus:" label with the short syllabus (english) text.
$shortSyllabus = @"
1)Heat-Exchange Equipment;
2)Tubular-Type Exchangers;
3)Plate-Type Exchangers;
4)Heat Transfer Systems Involving Phase Change;
5)Psychometrics;
6)Humidification and dehumidification operations
"@
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = $shortSyllabus
$ws.Range("C14").Value = $shortSyllabus
$ws.Rows.Item(14).RowHeight = 60

# Row 15: "Programa:" label, value becomes the first professor entry again.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = $teacher1
$ws.Range("C15").Value = $teacher1
$ws.Rows.Item(15).RowHeight = 120

# Row 16: "Syllabus:" label with the full (english) syllabus text.
$syllabus = @"
1)Heat Exchangers: general concepts and types of heat exchangers;
2)Tubular heat exchangers: calculations in a bitubular heat exchanger; LMTD method; Shell and tube heat exchangers; Correlations to determine the heat transfer coefficients in Shell and tube heat exchangers; Estimative of coefficients; NTU method;
3)Plate heat exchangers: calculations and comparison with tubular exchangers;
4)Heat exchange systems with phase change: evaporators, condensers, reboilers and boilers; crystallization;
5)Psychrometry: concepts involved and the use of psychrometric chart;
6)Humidification and dehumidification operations; Cooling towers and Drying
"@
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = $syllabus
$ws.Range("C16").Value = $syllabus
$ws.Rows.Item(16).RowHeight = 120

# Row 17: "Avaliacao:" label only, no values.
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Value = $null
$ws.Range("C17").Value = $null
$ws.Rows.Item(17).AutoFit()

# Row 18: "Metodo:" label, value becomes the second professor entry.
$teacher2 = "5817372 - Simone de Fátima Medeiros Sampaio"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = $teacher2
$ws.Range("C18").Value = $teacher2
$ws.Rows.Item(18).RowHeight = 60

# Row 19: "Criterio:" label, value becomes the exam-application text.
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aplicação de 2 provas, P1 e P2."
$ws.Range("C19").Value = "Aplicação de 2 provas, P1 e P2."
$ws.Rows.Item(19).RowHeight = 60

# Row 20: "Norma de recuperacao:" label, value becomes the grading criteria text.
$criterio = "A média do período será MP = (P1+P2)/2. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham frequência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou frequência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham frequência mínima de 70% serão submetidos ao período de recuperação (regimental)."
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio
$ws.Rows.Item(20).RowHeight = 60

# Row 21: "Bibliografia:" label, value becomes the recovery-norm text.
$normaRecuperacao = "A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação."
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = $normaRecuperacao
$ws.Range("C21").Value = $normaRecuperacao
$ws.Rows.Item(21).RowHeight = 120

# Row 22: "Requisitos:" label only, no values.
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Value = $null
$ws.Range("C22").Value = $null
$ws.Rows.Item(22).AutoFit()

# Row 23: no label, value becomes the weak-requirement text (trailing newline kept).
$requisito = "LOQ4084 -  Fenômenos de Transporte II  (Requisito fraco)`n"
$ws.Range("A23").Value = $null
$ws.Range("B23").Value = $requisito
$ws.Range("C23").Value = $requisito
$ws.Rows.Item(23).RowHeight = 30
